# Update countries & provincias Spain
# Refresh the "last updated" timestamp and the COVID-19 stats snapshot,
# which causes a handful of countries to change rank (row) order because
# the sheet is kept sorted by total cases (column B) descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "datos actualizados" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 04:01"

# Helper to write a full data row: País, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos críticos, Muertes hoy, Muertes
function Set-Row {
    param($row, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes)
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Rows 41-43: Bolivia overtakes Portugal and Singapur in total cases
Set-Row 41 "Bolivia"  47200 1635 14333 31113 0 52 1754
Set-Row 42 "Portugal" 46221 0    30655 13912 0 0  1654
Set-Row 43 "Singapur" 45783 0    42026 3731  0 0  26

# Row 67: Corea del Sur stats refresh (no reordering)
Set-Row 67 "Corea del Sur" 13417 44 12178 950 0 1 289

# Rows 74-75: Australia overtakes Kenia in total cases
Set-Row 74 "Australia" 9796 243 7727 1961 0 1 108
Set-Row 75 "Kenia"     9726 0   2832 6710 0 0 184

# Row 126: Nueva Zelanda stats refresh (no reordering)
Set-Row 126 "Nueva Zelanda" 1544 1 1497 25 0 0 22

# Rows 177-178: Camboya overtakes Bermudas in total cases
Set-Row 177 "Camboya"  156 15 133 23 0 0 0
Set-Row 178 "Bermudas" 150 0  137 4  0 0 9
